# Adding the Intraday Model for Imperial
# Shift the whole day of data from 19.06.2024 to 20.06.2024:
#   - Column A (Data): date/time serials move forward by exactly 1 day
#   - Column D (Lookup): text labels "19.06.2024<n>" -> "20.06.2024<n>"
#   - Column C (Prediction): refreshed model output for this new day
#     (only some rows' predictions actually change)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Prediction (column C) values for the rows whose forecast changed.
$newPredictions = @{
    27 = 0.016
    28 = 0.025
    29 = 0.035
    30 = 0.052
    31 = 0.08400000000000001
    32 = 0.12
    33 = 0.163
    34 = 0.211
    35 = 0.262
    36 = 0.316
    37 = 0.371
    38 = 0.423
    39 = 0.432
    40 = 0.477
    41 = 0.5659999999999999
    42 = 0.612
    43 = 0.64
    44 = 0.677
    45 = 0.716
    46 = 0.749
    47 = 0.776
    48 = 0.804
    49 = 0.824
    50 = 0.839
    51 = 0.856
    52 = 0.867
    53 = 0.871
    54 = 0.877
    55 = 0.873
    56 = 0.866
    57 = 0.875
    58 = 0.867
    59 = 0.847
    60 = 0.827
    73 = 0.414
    74 = 0.368
    75 = 0.313
    76 = 0.285
    77 = 0.222
    78 = 0.173
    79 = 0.146
    80 = 0.109
    81 = 0.09
    82 = 0.08
    83 = 0.061
    84 = 0.047
    85 = 0.034
}

for ($r = 2; $r -le 96; $r++) {
    # Column A: push the timestamp forward by one full day (serial + 1).
    $curDate = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $curDate + 1

    # Column C: overwrite prediction where the new day's model output differs.
    if ($newPredictions.ContainsKey($r)) {
        $ws.Cells.Item($r, 3).Value = $newPredictions[$r]
    }

    # Column D: rewrite the lookup label's date prefix from 19.06.2024 to 20.06.2024.
    $curLookup = $ws.Cells.Item($r, 4).Value2
    $newLookup = $curLookup -replace '^19\.06\.2024', '20.06.2024'
    $ws.Cells.Item($r, 4).Value = $newLookup
}
